$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add two new header values in P1 and Q1 ---
# Copy the format of the existing last header cell (O1, style index 1: bold/border/centered)
# onto the two new header cells, then set their values.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows 2..25 ---
# Columns I, K, M, O swap their values (1<->2) and two new columns P, Q are
# appended with value 2 (unstyled, like the other data columns B..O).
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P: new column = 2
    $ws.Cells.Item($r, 17).Value = 2   # Q: new column = 2
}
